$wb = $excel.ActiveWorkbook

# --- Sheet: Metadata ---
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("A2").Value = "30 Oct 2025, 10:03 AM"

# --- Sheet: 1 Month Performance ---
$wsPerf = $wb.Worksheets.Item("1 Month Performance")
$wsPerf.Range("C2").Value = 115.7484
$wsPerf.Range("C4").Value = 82.0808
$wsPerf.Range("C7").Value = 53.9642
$wsPerf.Range("C8").Value = 50.8554
$wsPerf.Range("C9").Value = 49.9287
$wsPerf.Range("C10").Value = 40.8226
$wsPerf.Range("C12").Value = 40.1101
$wsPerf.Range("C13").Value = 39.2925
$wsPerf.Range("C15").Value = 37.9079
$wsPerf.Range("C16").Value = 37.6534
$wsPerf.Range("C17").Value = 36.7527
$wsPerf.Range("C19").Value = 36.4424
$wsPerf.Range("C20").Value = 34.6996
$wsPerf.Range("C21").Value = 33.8319
$wsPerf.Range("C22").Value = 33.7255
$wsPerf.Range("C23").Value = 33.241
$wsPerf.Range("B24").Value = "ORIENTTECH"
$wsPerf.Range("C24").Value = 30.1278
$wsPerf.Range("B25").Value = "TARACHAND"
$wsPerf.Range("C25").Value = 29.9605
$wsPerf.Range("C26").Value = 28.92
$wsPerf.Range("C27").Value = 28.8967
$wsPerf.Range("C29").Value = 27.9985
$wsPerf.Range("C30").Value = 27.5131
$wsPerf.Range("C31").Value = 27.2849
$wsPerf.Range("C33").Value = 26.892
$wsPerf.Range("C34").Value = 26.4886
$wsPerf.Range("C35").Value = 25.9927
$wsPerf.Range("C36").Value = 25.0928
$wsPerf.Range("C37").Value = 25.0259
$wsPerf.Range("C38").Value = 24.8384
$wsPerf.Range("C41").Value = 23.5582
$wsPerf.Range("C42").Value = 23.5166
$wsPerf.Range("B43").Value = "DCBBANK"
$wsPerf.Range("C43").Value = 23.1404
$wsPerf.Range("C44").Value = 23.1061
$wsPerf.Range("B45").Value = "INDIANB"
$wsPerf.Range("C45").Value = 22.9258
$wsPerf.Range("C46").Value = 21.6991
$wsPerf.Range("C47").Value = 21.6932
$wsPerf.Range("B48").Value = "SURYODAY"
$wsPerf.Range("C48").Value = 21.6728
$wsPerf.Range("B50").Value = "INDRAMEDCO"
$wsPerf.Range("C50").Value = 21.3137
$wsPerf.Range("B51").Value = "GUJTHEM"
$wsPerf.Range("C51").Value = 21.2197
$wsPerf.Range("B53").Value = "RBLBANK"
$wsPerf.Range("C53").Value = 20.6884
$wsPerf.Range("B54").Value = "MOLDTECH"
$wsPerf.Range("C54").Value = 20.6214
$wsPerf.Range("B55").Value = "SCI"
$wsPerf.Range("C55").Value = 20.5958
$wsPerf.Range("B56").Value = "BLUEDART"
$wsPerf.Range("C56").Value = 20.5757
$wsPerf.Range("B57").Value = "BHARATWIRE"
$wsPerf.Range("C57").Value = 20.4967
$wsPerf.Range("B58").Value = "KERNEX"
$wsPerf.Range("C58").Value = 20.3666
$wsPerf.Range("B59").Value = "MARINE"
$wsPerf.Range("C59").Value = 20.1693
$wsPerf.Range("B60").Value = "SHRIRAMFIN"
$wsPerf.Range("C60").Value = 20.1343
$wsPerf.Range("B61").Value = "SAGILITY"
$wsPerf.Range("C61").Value = 20.076
$wsPerf.Range("B62").Value = "PRECWIRE"
$wsPerf.Range("C62").Value = 19.9568
$wsPerf.Range("C63").Value = 19.7354
$wsPerf.Range("B64").Value = "IIFL"
$wsPerf.Range("C64").Value = 19.7
$wsPerf.Range("C65").Value = 19.6407
$wsPerf.Range("B66").Value = "FEDERALBNK"
$wsPerf.Range("C66").Value = 19.5713
$wsPerf.Range("B67").Value = "THOMASCOTT"
$wsPerf.Range("C67").Value = 19.3679
$wsPerf.Range("C68").Value = 19.3513
$wsPerf.Range("C69").Value = 19.2461
$wsPerf.Range("B71").Value = "ETHOSLTD"
$wsPerf.Range("C71").Value = 19.0718
$wsPerf.Range("B72").Value = "UNIPARTS"
$wsPerf.Range("C72").Value = 18.9528
$wsPerf.Range("C73").Value = 18.6571
$wsPerf.Range("C74").Value = 18.2415
$wsPerf.Range("C75").Value = 18.1213
$wsPerf.Range("C76").Value = 17.7335

# --- Sheet: distance from Dma50 ---
$wsDma = $wb.Worksheets.Item("distance from Dma50")
$wsDma.Range("C2").Value = 9.4627
$wsDma.Range("C3").Value = 7.2767
$wsDma.Range("C4").Value = 5.7289
$wsDma.Range("C5").Value = 5.0341
$wsDma.Range("C6").Value = 4.9642
$wsDma.Range("C7").Value = 4.9408
$wsDma.Range("C8").Value = 4.3268
$wsDma.Range("C9").Value = 4.2614
$wsDma.Range("C10").Value = 3.511
$wsDma.Range("C11").Value = 3.4391
$wsDma.Range("C12").Value = 3.2361
$wsDma.Range("C13").Value = 3.108
$wsDma.Range("C14").Value = 3.0657
$wsDma.Range("C15").Value = 2.9553
$wsDma.Range("C16").Value = 2.9148
$wsDma.Range("C17").Value = 2.7352
$wsDma.Range("C18").Value = 2.5623
$wsDma.Range("C19").Value = 2.2211
$wsDma.Range("C20").Value = 2.1339
$wsDma.Range("C21").Value = 2.0913
$wsDma.Range("B22").Value = "NIFTYCONSUMPTION"
$wsDma.Range("C22").Value = 1.3281
$wsDma.Range("B23").Value = "CNXIT"
$wsDma.Range("C23").Value = 1.3179
$wsDma.Range("C24").Value = 1.0346
$wsDma.Range("C25").Value = 0.8193
$wsDma.Range("C26").Value = 0.8091
$wsDma.Range("C27").Value = 0.7066
$wsDma.Range("C28").Value = 0.1647
$wsDma.Range("C29").Value = -0.2117
$wsDma.Range("C30").Value = -1.9943
